$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Resultado" column (E) text for rows 5, 6, 8, 9 to include
# the parameter name in the "Argumento inválido" message, while row 7
# ("Não é um cliente") is untouched.
$ws.Range("E5").Value = "Argumento inválido  (Parameter 'nome')"
$ws.Range("E6").Value = "Argumento inválido  (Parameter 'nome')"
$ws.Range("E8").Value = "Argumento inválido  (Parameter 'idade')"
$ws.Range("E9").Value = "Argumento inválido  (Parameter 'idade')"

# Update the selection to match the new active cell/range in the sheet view.
$ws.Range("F5:F9").Select()
